$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper idea: this runtime auto-merges adjacent runs that end up with
# identical run properties (rPr). Word's real Find/Replace already does
# that (and conveniently also drops now-redundant <w:proofErr/> markers
# that used to bracket the individual spell-check-flagged words), so a
# single Find.Execute over the old/joined text is the cleanest way to
# fold a run-per-word paragraph into the merged text the diff wants.
#
# Where the diff instead wants a paragraph kept as *several* runs (all
# sharing the same rPr), a plain text edit collapses them to one run,
# so we nudge a sub-range's Bold property on/off right after - that
# forces the engine to keep that sub-range as its own run without
# altering the visible formatting (Bold ends up back at its original
# value).
# ---------------------------------------------------------------------

# 1) Footer byline: "Rasmus Tilljander - rati10@student.bth.se"
#    (collapse the spell-checked "Rasmus"/"Tilljander" runs + proofErr
#    markers into a single run)
$old1 = "Rasmus Tilljander - rati10@student.bth.se"
$r1 = $d.Content
$r1.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# 2) Work Summary: append a new sentence as its own trailing run
$p2 = $d.Paragraphs.Item(31)
$origText2 = "The whole week has gone to programming and some minor research into different programming solutions."
$addition2 = " We also wrote a summary of the research trip taken by some members of the team."
$pStart2 = $p2.Range.Start
$p2.Range.InsertAfter($addition2)
$newRun2 = $d.Range($pStart2 + $origText2.Length, $pStart2 + $origText2.Length + $addition2.Length)
$newRun2.Bold = 1
$newRun2.Bold = 0

# 3) Major design decisions: collapse the "SuperCandy" run-split + proofErr
$old3 = "We decided to take away the SuperCandy class and instead give the Candy class a type variable. The reason for this was that having two so similar classes seemed unnecessary. "
$r3 = $d.Content
$r3.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null

# 4) WBS changes: insert a new body paragraph right after the heading
$r4 = $d.Content
$r4.Find.Execute("WBS changes", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$headingPara4 = $r4.Paragraphs.Item(1)
$headingIdx4 = $headingPara4.Index
$headingPara4.Range.InsertParagraphAfter()
$newPara4 = $d.Paragraphs.Item($headingIdx4 + 1)
$newPara4.Range.Text = "No changes have been made to the WBS this week."
$newPara4.Range.LanguageID = "en-US"
$newPara4.Format.Style = "Normal"

# 5) Issues, problems and risks: reword + split into the diff's run layout
#    + append the new trailing sentence
$old5 = "currently we have a 50% processor usage even if we only run an empty messageloop in the main."
$new5 = "currently we have an alarmingly high processor usage even if we only run an empty messageloop in the main. If this problem is not resolved quickly it could stop the production almost completely since testing to see if functions work is made virtually impossible."
$r5 = $d.Content
$r5.Find.Execute($old5, $false, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

$r5b = $d.Content
$r5b.Find.Execute("currently we have an alarmingly high processor", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$para5 = $r5b.Paragraphs.Item(1)
$pStart5 = $para5.Range.Start

$seg1 = "currently we have a"
$seg2 = "n"
$seg3 = " "
$seg4 = "alarmingly high"
$seg5 = " processor usage even if we only run an empty messageloop in the main."
$seg6 = " If this problem is not resolved quickly it could stop the production almost completely since testing to see if functions work is made virtually impossible."
$lens5 = @($seg1.Length, $seg2.Length, $seg3.Length, $seg4.Length, $seg5.Length, $seg6.Length)

$cursor5 = $pStart5
$positions5 = @($cursor5)
foreach ($len in $lens5) {
    $cursor5 = $cursor5 + $len
    $positions5 += $cursor5
}
for ($i = 1; $i -lt $positions5.Length - 1; $i++) {
    $segRange = $d.Range($positions5[$i], $positions5[$i + 1])
    $segRange.Bold = 1
    $segRange.Bold = 0
}

# 6) Current status: collapse the ScreenHandler/PlayerInput/ResourceHandler/
#    GameTimer run-split + proofErr, but keep the diff's 3-run layout
$seg6a = "The Camera and the HUD is complete. The ScreenHandler is also finished for basic functions."
$seg6b = " The PlayerInput class is finished and so is the ResourceHandler."
$seg6c = " The GameTimer is finished."
$old6 = $seg6a + $seg6b + $seg6c
$r6 = $d.Content
$r6.Find.Execute($old6, $false, $false, $false, $false, $false, $true, 1, $false, $old6, 2) | Out-Null

$r6b = $d.Content
$r6b.Find.Execute($seg6a, $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$para6 = $r6b.Paragraphs.Item(1)
$pStart6 = $para6.Range.Start

$lens6 = @($seg6a.Length, $seg6b.Length, $seg6c.Length)
$cursor6 = $pStart6
$positions6 = @($cursor6)
foreach ($len in $lens6) {
    $cursor6 = $cursor6 + $len
    $positions6 += $cursor6
}
for ($i = 1; $i -lt $positions6.Length - 1; $i++) {
    $segRange = $d.Range($positions6[$i], $positions6[$i + 1])
    $segRange.Bold = 1
    $segRange.Bold = 0
}

# 7) Planned work: fill in the (previously empty) last paragraph
$p7 = $d.Paragraphs.Last
$p7.Range.Text = "Continued programming as well as starting work on the actual demo to be handed in."
$p7.Range.LanguageID = "en-US"

Write-Output "edit complete"
